# Update cryptos list prices (column D) and 1h volume/change percentages (column E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep these columns as plain text so numeric-looking strings (e.g. "112.60")
# are not silently coerced to numbers (which would drop trailing zeros).
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.691.14"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "1.855.76"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("D4").Value = "1.015"
$ws.Range("E4").Value = "  +1.12%  "
$ws.Range("D5").Value = "332.95"
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("D6").Value = "1.012"
$ws.Range("E6").Value = "  +0.87%  "
$ws.Range("D7").Value = "0.4681"
$ws.Range("E7").Value = "  -0.59%  "
$ws.Range("D8").Value = "0.3887"
$ws.Range("E8").Value = "  -1.31%  "
$ws.Range("D9").Value = "46.65"
$ws.Range("E9").Value = "  -2.63%  "
$ws.Range("D10").Value = "0.07946"
$ws.Range("E10").Value = "  -1.06%  "
$ws.Range("D11").Value = "1.001"
$ws.Range("E11").Value = "  -2.50%  "
$ws.Range("D12").Value = "21.51"
$ws.Range("E12").Value = "  -2.42%  "
$ws.Range("D13").Value = "1.872.09"
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("D14").Value = "5.988"
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").Value = "7.122"
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("D16").Value = "1.017"
$ws.Range("E16").Value = "  +0.93%  "
$ws.Range("D17").Value = "87.87"
$ws.Range("E17").Value = "  +1.14%  "
$ws.Range("D18").Value = "0.06708"
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("E19").Value = "  -0.73%  "
$ws.Range("D20").Value = "16.85"
$ws.Range("E20").Value = "  -1.69%  "
$ws.Range("E21").Value = "  +0.76%  "
$ws.Range("D22").Value = "27.690.06"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "5.448"
$ws.Range("E23").Value = "  -1.17%  "
$ws.Range("D24").Value = "10.89"
$ws.Range("E24").Value = "  -0.62%  "
$ws.Range("D25").Value = "2.318"
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("D26").Value = "2.090.33"
$ws.Range("E26").Value = "  -0.98%  "
$ws.Range("D27").Value = "158.49"
$ws.Range("E27").Value = "  -0.28%  "
$ws.Range("D28").Value = "19.63"
$ws.Range("E28").Value = "  -2.40%  "
$ws.Range("D29").Value = "2.083"
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("D30").Value = "5.383"
$ws.Range("E30").Value = "  -3.19%  "
$ws.Range("D31").Value = "120.65"
$ws.Range("E31").Value = "  -1.28%  "
$ws.Range("D32").Value = "0.9624"
$ws.Range("E32").Value = "  -1.07%  "
$ws.Range("D33").Value = "0.09424"
$ws.Range("E33").Value = "  -1.04%  "
$ws.Range("D34").Value = "3.643"
$ws.Range("E34").Value = "  +1.37%  "
$ws.Range("D35").Value = "5.291"
$ws.Range("E35").Value = "  -0.64%  "
$ws.Range("D36").Value = "1.329"
$ws.Range("E36").Value = "  -8.10%  "
$ws.Range("D37").Value = "0.06028"
$ws.Range("E37").Value = "  -1.20%  "
$ws.Range("D38").Value = "0.02208"
$ws.Range("E38").Value = "  -1.90%  "
$ws.Range("D39").Value = "1.208"
$ws.Range("E39").Value = "  -1.21%  "
$ws.Range("D40").Value = "8.138"
$ws.Range("E40").Value = "  +0.26%  "
$ws.Range("E41").Value = "  +0.80%  "
$ws.Range("D42").Value = "0.5902"
$ws.Range("E42").Value = "  -1.89%  "
$ws.Range("D43").Value = "0.1879"
$ws.Range("E43").Value = "  -1.01%  "
$ws.Range("D44").Value = "10.15"
$ws.Range("E44").Value = "  -0.86%  "
$ws.Range("D45").Value = "1.259"
$ws.Range("E45").Value = "  -0.51%  "
$ws.Range("D46").Value = "0.5598"
$ws.Range("E46").Value = "  -1.62%  "
$ws.Range("D47").Value = "12.04"
$ws.Range("E47").Value = "  -1.36%  "
$ws.Range("D48").Value = "1.909"
$ws.Range("E48").Value = "  -1.64%  "
$ws.Range("D49").Value = "3.295"
$ws.Range("E49").Value = "  -2.53%  "
$ws.Range("D50").Value = "0.06759"
$ws.Range("E50").Value = "  -1.91%  "
$ws.Range("D51").Value = "112.60"
